$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": the handoff step for the remaining
# "low" priority items (rows 4-7, the 1453bc42/384f40f2/774fd473/a3f22f1f
# files) just completed, so:
#   - their Priority (column E) flips from "low" to "ht"
#   - their Latest Handoff Datetime (column H) is refreshed to the
#     new handoff timestamp
# for both the "zh-cn" and "de-de" target-language sheets.

$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-25 10:32:53"
}

$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-25 10:33:01"
}

# The "Overview" sheet's "Latest HO Xliff Generate Date" column (G) for
# those same rows mirrors the de-de handoff timestamp.
$overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-25 10:33:01"
}
